# Slide 2, title shape ("제목 3") currently reads "실습 내용".
# Target: split into two runs -> "실습 내용 " (Korean) + "(Forking workflow)" (English),
# with an explicit trailing endParaRPr.

$p  = $ppt.ActivePresentation
$s  = $p.Slides.Item(2)
$sh = $s.Shapes.Item(2)

# Start from just the English portion so we can stamp its language while it is
# still the (only / first) run in the paragraph.
$tr = $sh.TextFrame.TextRange
$tr.Text = "(Forking workflow)"

$tr2 = $sh.TextFrame.TextRange
$tr2.LanguageID = "en-US"

# Prepend the Korean portion (inherits the run properties at this point); this
# becomes the new first run, so re-stamping the language now targets it.
$tr3 = $sh.TextFrame.TextRange
$null = $tr3.InsertBefore("실습 내용 ")

$tr4 = $sh.TextFrame.TextRange
$tr4.LanguageID = "ko-KR"
